$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 162, shifting existing rows 162..303 down to 163..304.
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new weekly price observation.
$ws.Range("A162").Value = 3
$ws.Range("B162").Value = "Femacal de La Calera"
$ws.Range("C162").Value = "Coquimbo"
$ws.Range("D162").Value = 44669
$ws.Range("E162").Value = 5
$ws.Range("F162").Value = 100112039
$ws.Range("G162").Value = "Ciboulette"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 85
$ws.Range("K162").Value = 2000
$ws.Range("L162").Value = 2000
$ws.Range("M162").Value = 2000
$ws.Range("N162").Value = "$/docena de atados"
$ws.Range("O162").Value = "Provincia de Quillota"
$ws.Range("P162").Value = 667
$ws.Range("Q162").Value = 3
$ws.Range("R162").Value = "Hortaliza"

# Match the date cell's number format style used by the rest of the column (style index 2).
$ws.Range("D162").NumberFormat = $ws.Range("D163").NumberFormat
